# Update "want to go" counts (column F) on several rows across sheets
# to reflect a newer data snapshot, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 15086
$ws1.Cells.Item(3, 6).Value = 19255
$ws1.Cells.Item(5, 6).Value = 150
$ws1.Cells.Item(14, 6).Value = 179
$ws1.Cells.Item(15, 6).Value = 233
$ws1.Cells.Item(17, 6).Value = 1489
$ws1.Cells.Item(20, 6).Value = 100
$ws1.Cells.Item(21, 6).Value = 240
$ws1.Cells.Item(22, 6).Value = 8057
$ws1.Cells.Item(27, 6).Value = 1255
$ws1.Cells.Item(28, 6).Value = 5
$ws1.Cells.Item(30, 6).Value = 6091
$ws1.Cells.Item(31, 6).Value = 120
$ws1.Cells.Item(35, 6).Value = 293
$ws1.Cells.Item(36, 6).Value = 5500
$ws1.Cells.Item(37, 6).Value = 1006
$ws1.Cells.Item(40, 6).Value = 52

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(3, 6).Value = 20

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 15086
$ws4.Cells.Item(3, 6).Value = 19255
$ws4.Cells.Item(5, 6).Value = 150
$ws4.Cells.Item(14, 6).Value = 179
$ws4.Cells.Item(15, 6).Value = 233
$ws4.Cells.Item(17, 6).Value = 1489
$ws4.Cells.Item(21, 6).Value = 100
$ws4.Cells.Item(22, 6).Value = 240
$ws4.Cells.Item(23, 6).Value = 8057
$ws4.Cells.Item(28, 6).Value = 1255
$ws4.Cells.Item(29, 6).Value = 5
$ws4.Cells.Item(31, 6).Value = 20
$ws4.Cells.Item(33, 6).Value = 6091
$ws4.Cells.Item(34, 6).Value = 120
$ws4.Cells.Item(38, 6).Value = 293
$ws4.Cells.Item(39, 6).Value = 5500
$ws4.Cells.Item(40, 6).Value = 1006
$ws4.Cells.Item(43, 6).Value = 52

$wb.Save()
